# Updated symbol list on Fri Feb 17 17:23:41 UTC 2023 with GitHub Actions
#
# The published crypto snapshot refreshed its prices/volumes; a few rows of
# the ranking also reshuffled (several coins swapped ranks), so the B/C/D/E
# cells for some rows now hold what used to be a neighbouring row's data,
# with freshly refreshed Price/Volume(1h) figures.
#
# Every value below is written as literal TEXT (prices and percentages are
# stored as strings in this sheet already, not as numbers). Setting
# NumberFormat to "@" before the assignment keeps numeric-looking text (e.g.
# "310.47" or "-3.23%") from being parsed into a number/date by Excel, and
# resetting the Style back to "Normal" afterwards means we do not leave a
# stray "Text"-formatted style behind on cells that were unstyled before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '310.47'
$cell.Style = "Normal"

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = '-3.23%'
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '54.35'
$cell.Style = "Normal"

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = '10.12%'
$cell.Style = "Normal"

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '5.105'
$cell.Style = "Normal"

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = '-4.32%'
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '0.07902'
$cell.Style = "Normal"

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = '-1.65%'
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '4.541'
$cell.Style = "Normal"

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = '-1.50%'
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '1.391'
$cell.Style = "Normal"

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = '1.74%'
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '1.704'
$cell.Style = "Normal"

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = '3.99%'
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.1245'
$cell.Style = "Normal"

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = '-1.71%'
$cell.Style = "Normal"

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.2021'
$cell.Style = "Normal"

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = '2.96%'
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.04744'
$cell.Style = "Normal"

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = '0.54%'
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.09430'
$cell.Style = "Normal"

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = '-3.49%'
$cell.Style = "Normal"

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.1049'
$cell.Style = "Normal"

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = '0.22%'
$cell.Style = "Normal"

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.001273'
$cell.Style = "Normal"

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = '-3.86%'
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.005679'
$cell.Style = "Normal"

$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = '-3.15%'
$cell.Style = "Normal"

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = '2,011.75%'
$cell.Style = "Normal"

$cell = $ws.Range("B17")
$cell.NumberFormat = "@"
$cell.Value = 'LEO'
$cell.Style = "Normal"

$cell = $ws.Range("C17")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '3.326'
$cell.Style = "Normal"

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = '-0.64%'
$cell.Style = "Normal"

$cell = $ws.Range("B18")
$cell.NumberFormat = "@"
$cell.Value = 'BTSEToken'
$cell.Style = "Normal"

$cell = $ws.Range("C18")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$cell.Style = "Normal"

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '2.417'
$cell.Style = "Normal"

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = '-1.45%'
$cell.Style = "Normal"

$cell = $ws.Range("B19")
$cell.NumberFormat = "@"
$cell.Value = 'BitpandaEcosystemToken'
$cell.Style = "Normal"

$cell = $ws.Range("C19")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.3424'
$cell.Style = "Normal"

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = '-2.31%'
$cell.Style = "Normal"

$cell = $ws.Range("B20")
$cell.NumberFormat = "@"
$cell.Value = 'MCDex'
$cell.Style = "Normal"

$cell = $ws.Range("C20")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '8.330'
$cell.Style = "Normal"

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = '3.25%'
$cell.Style = "Normal"

$cell = $ws.Range("B21")
$cell.NumberFormat = "@"
$cell.Value = 'ProBitToken'
$cell.Style = "Normal"

$cell = $ws.Range("C21")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '0.1361'
$cell.Style = "Normal"

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = '0.04%'
$cell.Style = "Normal"

$cell = $ws.Range("B22")
$cell.NumberFormat = "@"
$cell.Value = 'ZBToken'
$cell.Style = "Normal"

$cell = $ws.Range("C22")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.2900'
$cell.Style = "Normal"

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = '-6.26%'
$cell.Style = "Normal"

$cell = $ws.Range("B23")
$cell.NumberFormat = "@"
$cell.Value = 'CoinExToken'
$cell.Style = "Normal"

$cell = $ws.Range("C23")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.04161'
$cell.Style = "Normal"

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = '-0.90%'
$cell.Style = "Normal"

$cell = $ws.Range("B24")
$cell.NumberFormat = "@"
$cell.Value = 'BitKan'
$cell.Style = "Normal"

$cell = $ws.Range("C24")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '0.001254'
$cell.Style = "Normal"

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = '-4.77%'
$cell.Style = "Normal"

$cell = $ws.Range("B25")
$cell.NumberFormat = "@"
$cell.Value = 'HotbitToken'
$cell.Style = "Normal"

$cell = $ws.Range("C25")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$cell.Style = "Normal"

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.003990'
$cell.Style = "Normal"

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = '-7.91%'
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '0.0001343'
$cell.Style = "Normal"

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = '-0.75%'
$cell.Style = "Normal"

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.02620'
$cell.Style = "Normal"

$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = '-4.02%'
$cell.Style = "Normal"

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.05956'
$cell.Style = "Normal"

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = '-2.26%'
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.009868'
$cell.Style = "Normal"

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = '-9.34%'
$cell.Style = "Normal"

$cell = $ws.Range("B41")
$cell.NumberFormat = "@"
$cell.Value = 'BKEXToken'
$cell.Style = "Normal"

$cell = $ws.Range("C41")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.1748'
$cell.Style = "Normal"

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = '19.46%'
$cell.Style = "Normal"

$cell = $ws.Range("B42")
$cell.NumberFormat = "@"
$cell.Value = 'KickToken'
$cell.Style = "Normal"

$cell = $ws.Range("C42")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.007995'
$cell.Style = "Normal"

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = '-0.81%'
$cell.Style = "Normal"

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.008186'
$cell.Style = "Normal"

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = '3.52%'
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.008317'
$cell.Style = "Normal"

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = '-4.28%'
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.3388'
$cell.Style = "Normal"

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = '-2.85%'
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.00007243'
$cell.Style = "Normal"

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = '5.58%'
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.00000000746'
$cell.Style = "Normal"

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = '-0.72%'
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.09442'
$cell.Style = "Normal"

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = '59.24%'
$cell.Style = "Normal"

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.002610'
$cell.Style = "Normal"

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = '-34.86%'
$cell.Style = "Normal"

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.00002089'
$cell.Style = "Normal"

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = '-0.72%'
$cell.Style = "Normal"

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.0001989'
$cell.Style = "Normal"

$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = '-0.72%'
$cell.Style = "Normal"
